# edit.ps1 - applies the two text edits described by the diff:
#
#  1) Slide 10 ("Under [Project]:" body placeholder), paragraph:
#       "This should be consistent with the base file name used previously"
#     -> split into three runs:
#       "This should be consistent with the base file "
#       "name and project folder named "
#       "previously"
#
#  2) Slide 4 ("Create and name a project folder..." body placeholder),
#     paragraph made of two runs:
#       "Ensure the Cfturbo Design File is also in the "
#       "project folder, both with the same name!"
#     -> merged into a single run containing the full sentence.
#
# Note: PowerPoint's TextRange.Text getter returns each paragraph's text
# with a trailing CR (chr 13), so comparisons are trimmed before use.

$p = $ppt.ActivePresentation

# --- Edit 1: Slide 10 -----------------------------------------------------
$slide10 = $p.Slides.Item(10)
$shape10 = $slide10.Shapes.Item(1)
$tr10 = $shape10.TextFrame.TextRange

$oldFragment = "name used "
$newFragment = "name and project folder named "

for ($i = 1; $i -le $tr10.Paragraphs().Count; $i++) {
    $para = $tr10.Paragraphs($i)
    $idx = $para.Text.IndexOf($oldFragment)
    if ($idx -ge 0) {
        $sub = $para.Characters($idx + 1, $oldFragment.Length)
        $sub.Text = $newFragment
        break
    }
}

# --- Edit 2: Slide 4 -------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$shape4 = $slide4.Shapes.Item(1)
$tr4 = $shape4.TextFrame.TextRange

$mergedText = "Ensure the Cfturbo Design File is also in the project folder, both with the same name!"

for ($i = 1; $i -le $tr4.Paragraphs().Count; $i++) {
    $para = $tr4.Paragraphs($i)
    $current = $para.Text.TrimEnd("`r", "`n")
    if ($current -eq $mergedText) {
        # Re-setting the paragraph text to its own concatenated value is
        # normally a no-op (text already matches) and the two pre-existing
        # runs would stay split, so force a real change first, then write
        # the final text back to let the engine coalesce into one run.
        $para.Text = "placeholder"
        $para.Text = $mergedText
        break
    }
}
